$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 716, pushing existing rows 716:798 down to 720:802
$ws.Rows("716:719").Insert()

# Fill the 4 new rows (716-719) with the new weekly data points.
# Columns A,B,C,E,F,G,H,I,J,Q,R,T are constant across this sheet's data rows.

# Row 716: Naranja / Lane Late / Primera
$ws.Range("A716").Value2 = 8
$ws.Range("B716").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C716").Value2 = "Coquimbo"
$ws.Range("D716").Value2 = 44826
$ws.Range("E716").Value2 = 4
$ws.Range("F716").Value2 = "Fruta"
$ws.Range("G716").Value2 = 100102
$ws.Range("H716").Value2 = "Cítricos"
$ws.Range("I716").Value2 = 100102005
$ws.Range("J716").Value2 = "Naranja"
$ws.Range("K716").Value2 = "Lane Late"
$ws.Range("L716").Value2 = "Primera"
$ws.Range("M716").Value2 = 22
$ws.Range("N716").Value2 = 85000
$ws.Range("O716").Value2 = 90000
$ws.Range("P716").Value2 = 87500
$ws.Range("Q716").Value2 = "$/bins (400 kilos)"
$ws.Range("R716").Value2 = "Provincia de Limarí"
$ws.Range("S716").Value2 = 219
$ws.Range("T716").Value2 = 400

# Row 717: Naranja / Lane Late / Segunda
$ws.Range("A717").Value2 = 8
$ws.Range("B717").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C717").Value2 = "Coquimbo"
$ws.Range("D717").Value2 = 44826
$ws.Range("E717").Value2 = 4
$ws.Range("F717").Value2 = "Fruta"
$ws.Range("G717").Value2 = 100102
$ws.Range("H717").Value2 = "Cítricos"
$ws.Range("I717").Value2 = 100102005
$ws.Range("J717").Value2 = "Naranja"
$ws.Range("K717").Value2 = "Lane Late"
$ws.Range("L717").Value2 = "Segunda"
$ws.Range("M717").Value2 = 16
$ws.Range("N717").Value2 = 55000
$ws.Range("O717").Value2 = 60000
$ws.Range("P717").Value2 = 57500
$ws.Range("Q717").Value2 = "$/bins (400 kilos)"
$ws.Range("R717").Value2 = "Provincia de Limarí"
$ws.Range("S717").Value2 = 144
$ws.Range("T717").Value2 = 400

# Row 718: Naranja / Navel Late / Primera
$ws.Range("A718").Value2 = 8
$ws.Range("B718").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C718").Value2 = "Coquimbo"
$ws.Range("D718").Value2 = 44826
$ws.Range("E718").Value2 = 4
$ws.Range("F718").Value2 = "Fruta"
$ws.Range("G718").Value2 = 100102
$ws.Range("H718").Value2 = "Cítricos"
$ws.Range("I718").Value2 = 100102005
$ws.Range("J718").Value2 = "Naranja"
$ws.Range("K718").Value2 = "Navel Late"
$ws.Range("L718").Value2 = "Primera"
$ws.Range("M718").Value2 = 24
$ws.Range("N718").Value2 = 85000
$ws.Range("O718").Value2 = 90000
$ws.Range("P718").Value2 = 87500
$ws.Range("Q718").Value2 = "$/bins (400 kilos)"
$ws.Range("R718").Value2 = "Provincia de Limarí"
$ws.Range("S718").Value2 = 219
$ws.Range("T718").Value2 = 400

# Row 719: Naranja / Navel Late / Segunda
$ws.Range("A719").Value2 = 8
$ws.Range("B719").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C719").Value2 = "Coquimbo"
$ws.Range("D719").Value2 = 44826
$ws.Range("E719").Value2 = 4
$ws.Range("F719").Value2 = "Fruta"
$ws.Range("G719").Value2 = 100102
$ws.Range("H719").Value2 = "Cítricos"
$ws.Range("I719").Value2 = 100102005
$ws.Range("J719").Value2 = "Naranja"
$ws.Range("K719").Value2 = "Navel Late"
$ws.Range("L719").Value2 = "Segunda"
$ws.Range("M719").Value2 = 16
$ws.Range("N719").Value2 = 55000
$ws.Range("O719").Value2 = 60000
$ws.Range("P719").Value2 = 57500
$ws.Range("Q719").Value2 = "$/bins (400 kilos)"
$ws.Range("R719").Value2 = "Provincia de Limarí"
$ws.Range("S719").Value2 = 144
$ws.Range("T719").Value2 = 400

# Ensure the date cells keep the expected date number format (style index 2 in the
# original file, numFmtId 165 "YYYY-MM-DD HH:MM:SS"). Insert() above already carries
# the formatting down from row 715, so explicitly mirror it to be safe.
$ws.Range("D716:D719").NumberFormat = $ws.Range("D715").NumberFormat

Write-Host ("Dimension after edit: " + $ws.UsedRange.Address())
